$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10, columns E (5) through T (20), matching the
# updated TPM-based recomputation described in the commit.
$data = @(
    @(3,1,0.331241,0.9937229999999999,0.2019242199214145,0.2019242199214145,1,0.3333333333333333,0.023286,0.069858,0.009310710475795457,0.009310710475795458,0.007713277925999999,0.069419501334,0.00188005794973914,0.00188005794973914),
    @(3,1,0.331241,0.9937229999999999,0.2019242199214145,0.2019242199214145,2,0.6666666666666666,0.1329193333333333,0.3987579999999999,0.05314667307834813,0.05314667307834814,0.04402833289266665,0.3962549960339999,0.01073160050276389,0.01073160050276389),
    @(3,1,0.331241,0.9937229999999999,0.2019242199214145,0.2019242199214145,3,1,2.344785333333334,7.034356000000001,0.9375426164458565,0.9375426164458565,0.7766890385986667,6.990201347388,0.1893125614689115,0.1893125614689115),
    @(3,1,0.9891043333333333,2.967313,0.6029571246591579,0.6029571246591579,1,0.3333333333333333,0.023286,0.069858,0.009310710475795457,0.009310710475795458,0.023032283506,0.207290551554,0.005613959217019528,0.005613959217019529),
    @(3,1,0.9891043333333333,2.967313,0.6029571246591579,0.6029571246591579,2,0.6666666666666666,0.1329193333333333,0.3987579999999999,0.05314667307834813,0.05314667307834814,0.1314710885837777,1.183239797254,0.03204516518452107,0.03204516518452107),
    @(3,1,0.9891043333333333,2.967313,0.6029571246591579,0.6029571246591579,3,1,2.344785333333334,7.034356000000001,0.9375426164458565,0.9375426164458565,2.319237333936444,20.873136005428,0.5652980002576173,0.5652980002576173),
    @(2,0.6666666666666666,0.320077,0.9602310000000001,0.1951186554194276,0.1951186554194276,1,0.3333333333333333,0.023286,0.069858,0.009310710475795457,0.009310710475795458,0.007453313022,0.06707981719800001,0.001816693309036788,0.001816693309036788),
    @(2,0.6666666666666666,0.320077,0.9602310000000001,0.1951186554194276,0.1951186554194276,2,0.6666666666666666,0.1329193333333333,0.3987579999999999,0.05314667307834813,0.05314667307834814,0.04254442145533333,0.382899793098,0.01036990739106318,0.01036990739106318),
    @(2,0.6666666666666666,0.320077,0.9602310000000001,0.1951186554194276,0.1951186554194276,3,1,2.344785333333334,7.034356000000001,0.9375426164458565,0.9375426164458565,0.7505118551373334,6.754606696236001,0.1829320547193276,0.1829320547193276)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 5   # column E = 5
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
